$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.430.36"
$ws.Range("E2").Value = "  +0.26%  "
# Row 3
$ws.Range("D3").Value = "1.868.40"
$ws.Range("E3").Value = "  -0.62%  "
# Row 4
$ws.Range("E4").Value = "  -0.16%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7068"
$ws.Range("E5").Value = "  -0.66%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.46"
$ws.Range("E6").Value = "  +0.39%  "
# Row 7
$ws.Range("E7").Value = "  -0.15%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07866"
$ws.Range("E8").Value = "  -1.90%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3133"
$ws.Range("E9").Value = "  -0.95%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.52"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07973"
$ws.Range("E11").Value = "  -4.08%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.211"
$ws.Range("E12").Value = "  -0.84%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.848.23"
$ws.Range("E13").Value = "  -1.80%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.38"
$ws.Range("E14").Value = "  -1.20%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7003"
$ws.Range("E15").Value = "  -2.10%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.529"
$ws.Range("E16").Value = "  +2.75%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008381"
$ws.Range("E17").Value = "  -2.02%  "
# Row 18
$ws.Range("D18").Value = "29.411.56"
$ws.Range("E18").Value = "  +0.21%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.20"
$ws.Range("E19").Value = "  +3.51%  "
# Row 20
$ws.Range("D20").Value = "2.125.52"
$ws.Range("E20").Value = "  -0.01%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -1.40%  "
# Row 22
$ws.Range("E22").Value = "  -0.12%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.625"
$ws.Range("E23").Value = "  -2.34%  "
# Row 24
$ws.Range("E24").Value = "  -0.24%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1553"
$ws.Range("E25").Value = "  -0.47%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.005"
$ws.Range("E26").Value = "  -0.89%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.31"
$ws.Range("E27").Value = "  -0.95%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("E28").Value = "  +0.77%  "
# Row 29
$ws.Range("E29").Value = "  -0.34%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.326"
$ws.Range("E30").Value = "  -2.21%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.265"
$ws.Range("E31").Value = "  -1.39%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +1.47%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05315"
$ws.Range("E33").Value = "  -1.38%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.894"
$ws.Range("E34").Value = "  -2.26%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7501"
$ws.Range("E35").Value = "  -2.55%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.173"
$ws.Range("E36").Value = "  -0.93%  "
# Row 37
$ws.Range("E37").Value = "  +1.03%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  -0.13%  "
# Row 39
$ws.Range("D39").Value = "1.276.51"
$ws.Range("E39").Value = "  +1.27%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.747"
$ws.Range("E40").Value = "  -0.23%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8940"
$ws.Range("E41").Value = "  -1.24%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.060"
$ws.Range("E42").Value = "  -6.74%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.74"
$ws.Range("E43").Value = "  -3.93%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.31"
$ws.Range("E44").Value = "  -4.05%  "
# Row 45
$ws.Range("E45").Value = "  -0.15%  "
# Row 46
$ws.Range("E46").Value = "  -3.98%  "
# Row 47
$ws.Range("D47").Value = "2.025.78"
$ws.Range("E47").Value = "  +0.33%  "
# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.793"
$ws.Range("E48").Value = "  -0.54%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.559"
$ws.Range("E49").Value = "  +1.06%  "
# Row 50
$ws.Range("E50").Value = "  -0.95%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4310"
$ws.Range("E51").Value = "  -1.56%  "
